$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 4, shifting rows 4-8 down to 5-9.
$ws.Rows.Item(4).Insert()

# Row 4: new "Zusatzbeitrag Krankenversicherung AG-Anteil in Prozent"
$ws.Range("A4").Value = "Zusatzbeitrag Krankenversicherung AG-Anteil in Prozent"
$ws.Range("B4").NumberFormat = "0.00"
$ws.Range("B4").Value = 0.99

# Row 5 (was row4): AN-Anteil
$ws.Range("A5").Value = "Zusatzbeitrag Krankenversicherung AN-Anteil in Prozent"
$ws.Range("B5").Value = 0.99

# Row 6 (was row5): Umlage U1 in Prozent
$ws.Range("A6").Value = "Umlage U1 in Prozent"
$ws.Range("B6").Value = 2.3

# Row 7 (was row6): Umlage U2 in Prozent
$ws.Range("A7").Value = "Umlage U2 in Prozent"
$ws.Range("B7").Value = 0.44

# Row 8 (was row7): Insolvenzgeldumlage
$ws.Range("A8").Value = "Insolvenzgeldumlage"
$ws.Range("B8").NumberFormat = "0.00"
$ws.Range("B8").Value = 0.06

# Row 9 (was row8, Eintragungsdatum / date) - update date value, keep text style
$ws.Range("A9").Value = "Eintragungsdatum"
$ws.Range("B9").Value = "01.01.2024"

$ws.Range("B5").Select()
